$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the "median levels of confidence" results paragraph.
# ---------------------------------------------------------------------------
$oldPara = "The median levels of confidence for organizing quantitative data using R and cleaning quantitative data using R increased from 2 (somewhat confident) to 3 (confident; p = 0.04) and 2 to 3 (p = 0.04), respectively. We also saw increases in the median level of confidence for using software to analyze quantitative data, visualize quantitative data,  and write a methods section, though these increases did not reach statistical significance. "
$newPara = "The median levels of confidence for organizing quantitative data using R (Task A), cleaning quantitative data using R (Task B), analyzing quantitative data using R (Task C), and visualizing quantitative data in R (Task D) each increased from 2 (somewhat confident) to 3 (confident; all p values < 0.05). We also saw significant increases in reported confidence in ability to write methods and results sections (Task E and Task F).  Median levels of confidence in the ability to share work in a reproducible way (Task G) and findings answers to questions (Task H) also increased, but these increases were not statistically significant. "

$d.Content.Find.Execute($oldPara, $true, $false, $false, $false, $false, $true, 1, $false, $newPara, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Table 1 updates (Pre-class rows for Tasks A, B, C, D, E, F, G): the
#    "4: Very confident" (or "1: Not confident at all" for Task F) N (%)
#    column and the paired p-value column both change.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Task A - Pre-class (row 4): col 6 = "4: Very confident", col 8 = p value
$t.Cell(4, 6).Range.Text = "1 (8.3)"
$t.Cell(4, 8).Range.Text = "0.02"

# Task B - Pre-class (row 6)
$t.Cell(6, 6).Range.Text = "1 (8.3)"
$t.Cell(6, 8).Range.Text = "0.02"

# Task C - Pre-class (row 8)
$t.Cell(8, 6).Range.Text = "1 (8.3)"
$t.Cell(8, 8).Range.Text = "0.03"

# Task D - Pre-class (row 10)
$t.Cell(10, 6).Range.Text = "1 (8.3)"
$t.Cell(10, 8).Range.Text = "0.03"

# Task E - Pre-class (row 12)
$t.Cell(12, 6).Range.Text = "1 (8.3)"
$t.Cell(12, 8).Range.Text = "0.03"

# Task F - Pre-class (row 14)
$t.Cell(14, 6).Range.Text = "2 (16.7)"
$t.Cell(14, 8).Range.Text = "0.04"

# Task G - Pre-class (row 16)
$t.Cell(16, 6).Range.Text = "1 (8.3)"
$t.Cell(16, 8).Range.Text = "0.11"

# ---------------------------------------------------------------------------
# 3. Resize the bar chart picture and clear its (now stale) auto-generated
#    alt text / description. Re-fetch the InlineShapes collection since the
#    earlier table/text edits can invalidate previously held handles.
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(1).LockAspectRatio = 0
$d.InlineShapes.Item(1).Width = 330.52228346456695
$d.InlineShapes.Item(1).Height = 220.34818897637794
$d.InlineShapes.Item(1).AlternativeText = ""
